$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "24"
$ws.Range("A20").Style = "Normal"
$ws.Range("B20").Value = "[BUG] <title>ahjkdahjsd"
$ws.Range("C20").Value = "open"
$ws.Range("D20").Value = "2025-03-26T06:55:02Z"
$ws.Range("E20").Value = "bug"
